# edit.ps1 - reproduce the authored change:
#   1. Slide 16's table switches from table style {C7D5181D-...} to
#      {A423CD0E-4F9A-41C8-A5E5-B2ACDDFB561E}.
#   2. The deck's theme colour palette changes from the "Integral" palette
#      to the stock "Office Theme" palette (dk1/lt1 are already identical
#      between the two palettes, so only the other ten slots move).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 ------------------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A423CD0E-4F9A-41C8-A5E5-B2ACDDFB561E}")
    }
}

# --- 2. Theme colour palette: Integral -> Office Theme --------------------
# Table mapping MsoThemeColorSchemeIndex -> RGB (dk1/lt1 unchanged: 000000 / FFFFFF).
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
